$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells for team record columns
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy header formatting (bold font + border) from an existing header cell
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill in the team record values for every data row (2-41)
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 29).Value = 87
    $ws.Cells.Item($r, 30).Value = 75
    $ws.Cells.Item($r, 31).Value = 0
}
